$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds daily data through row 238 (2021-04-26).
# Append six more days of data (through 2021-05-02), matching the
# "aggiornamento fino a 02/05" update. Column A keeps the same date
# style (s="2") as the preceding rows, so copy the formatting down
# from the last existing row before writing the new values.

$lastRow = 238
$newRows = @(
    @{ A = 44313; B = 0; C = 10; D = 312.5976867771179 },
    @{ A = 44314; B = 0; C = 8;  D = 250.0781494216943 },
    @{ A = 44315; B = 0; C = 6;  D = 187.5586120662707 },
    @{ A = 44316; B = 0; C = 6;  D = 187.5586120662707 },
    @{ A = 44317; B = 0; C = 5;  D = 156.2988433885589 },
    @{ A = 44318; B = 1; C = 1;  D = 31.25976867771178 }
)

$r = $lastRow
foreach ($row in $newRows) {
    $r = $r + 1
    $ws.Range("A$lastRow").Copy($ws.Range("A$r"))
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
}
